$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("A2").Value = "DELIGHT MILK"
$ws.Range("B2").Value = 1800

# Add new rows 3-5
$ws.Range("A3").Value = "NESTLE MILK"
$ws.Range("B3").Value = 1250

$ws.Range("A4").Value = "HST MILK "
$ws.Range("B4").Value = 250

$ws.Range("A5").Value = "MILMA"
$ws.Range("B5").Value = 600
